# Correção das notas do fórum para matc65 em 2021.2
# Zera todas as notas diárias de visualização (colunas B a J) para as linhas 2 a 50,
# mantendo a coluna A (matrícula) e o cabeçalho da linha 1 intactos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
